$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 1).Value = "7b1DtX5u"
$ws.Cells.Item(4, 2).Value = "trashboatsr"
$ws.Cells.Item(4, 3).Value = 1818
$ws.Cells.Item(4, 4).Value = 170
$ws.Cells.Item(4, 5).Value = "https://lichess.org/7b1DtX5u"
$ws.Cells.Item(4, 6).Value = 2886
$ws.Cells.Item(4, 7).Value = $true
$ws.Cells.Item(4, 8).Value = "newguy"

# Row 5
$ws.Cells.Item(5, 1).Value = "GHSOa063"
$ws.Cells.Item(5, 2).Value = "trashboatsr"
$ws.Cells.Item(5, 3).Value = 1818
$ws.Cells.Item(5, 4).Value = 140
$ws.Cells.Item(5, 5).Value = "https://lichess.org/GHSOa063"
$ws.Cells.Item(5, 6).Value = 2887
$ws.Cells.Item(5, 7).Value = $true
$ws.Cells.Item(5, 8).Value = "BLANK"

# Row 6
$ws.Cells.Item(6, 1).Value = "tXRhbtBp"
$ws.Cells.Item(6, 2).Value = "trashboatsr"
$ws.Cells.Item(6, 3).Value = 1818
$ws.Cells.Item(6, 4).Value = 180
$ws.Cells.Item(6, 5).Value = "https://lichess.org/tXRhbtBp"
$ws.Cells.Item(6, 6).Value = 2889
$ws.Cells.Item(6, 7).Value = $false
$ws.Cells.Item(6, 8).Value = "blank"

# Row 7
$ws.Cells.Item(7, 1).Value = "h3KmFj3w"
$ws.Cells.Item(7, 2).Value = "trashboatsr"
$ws.Cells.Item(7, 3).Value = 1818
$ws.Cells.Item(7, 4).Value = 210
$ws.Cells.Item(7, 5).Value = "https://lichess.org/h3KmFj3w"
$ws.Cells.Item(7, 6).Value = 2890
$ws.Cells.Item(7, 7).Value = $false
$ws.Cells.Item(7, 8).Value = "blank"

# Row 8
$ws.Cells.Item(8, 1).Value = "IJQleUwj"
$ws.Cells.Item(8, 2).Value = "trashboatsr"
$ws.Cells.Item(8, 3).Value = 1818
$ws.Cells.Item(8, 4).Value = 100
$ws.Cells.Item(8, 5).Value = "https://lichess.org/IJQleUwj"
$ws.Cells.Item(8, 6).Value = 2912
$ws.Cells.Item(8, 7).Value = $false
$ws.Cells.Item(8, 8).Value = "blank"
